$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("B2:F27")
$keyRange = $ws.Range("D2:D27")

$sortRange.Sort($keyRange, 1)

$ws.Range("D7").Select()
